# Weekly update: insert a new Coliflor price record for
# "Macroferia Regional de Talca" as the new first row of the data block
# (row 151), pushing the existing rows 151-178 down to 152-179.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151 (shifts old row 151.. down by one,
# and copies formatting - e.g. the date number format - from the
# surrounding rows automatically).
$ws.Rows(151).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value = 44543
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = 100112008
$ws.Range("G151").Value = "Coliflor"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 3000
$ws.Range("K151").Value = 900
$ws.Range("L151").Value = 900
$ws.Range("M151").Value = 900
$ws.Range("N151").Value = "$/unidad"
$ws.Range("O151").Value = "Región del Maule"
$ws.Range("P151").Value = 900
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = "Hortaliza"
